# "create contact and edit contact completed"
#
# The contact list previously had 2 rows (Anisree, Vaishak). Vaishak's
# contact record is edited (phone + pincode), Anisree's contact record is
# edited (phone, DOB, address, pincode) and a brand-new contact (Arjun) is
# created as row 4.
#
# Final layout:
#   Row2: VaishakT S | vaishak@gmail.com | 9876543210 | Male | 2001-04-10 | MarthandamTamil Nadu | 654321
#   Row3: AnisreeS S | anisree@gmail.com | 9876543210 | Male | 2000-07-12 | VenjaranmooduTVM      | 654321
#   Row4: ArjunK     | arjun@gmail.com   | 9638527410 | Male | 2001-04-16 | AnachalIdukki         | 123456

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Writes $value into $addr as a literal text cell, even when $value
    # looks like a date/number (e.g. "2001-04-10"), by entering it as a
    # string-literal formula and then collapsing the formula to its
    # computed value via a values-only paste. Plain Range.Value = "..."
    # would otherwise let Excel's auto-detection silently reinterpret a
    # date-shaped string as a date serial number.
    param($addr, $value)

    $escaped = $value.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

# ---- Row 2: now VaishakT S (moved up from row 3, phone + pincode edited)
Set-TextValue "A2" "VaishakT S"
Set-TextValue "B2" "vaishak@gmail.com"
$ws.Range("C2").Value = 9876543210
Set-TextValue "D2" "Male"
Set-TextValue "E2" "2001-04-10"
Set-TextValue "F2" "MarthandamTamil Nadu"
$ws.Range("G2").Value = 654321

# ---- Row 3: now AnisreeS S (moved down from row 2, several fields edited)
Set-TextValue "A3" "AnisreeS S"
Set-TextValue "B3" "anisree@gmail.com"
$ws.Range("C3").Value = 9876543210
Set-TextValue "D3" "Male"
Set-TextValue "E3" "2000-07-12"
Set-TextValue "F3" "VenjaranmooduTVM"
$ws.Range("G3").Value = 654321

# ---- Row 4: brand-new contact, ArjunK
Set-TextValue "A4" "ArjunK"
Set-TextValue "B4" "arjun@gmail.com"
$ws.Range("C4").Value = 9638527410
Set-TextValue "D4" "Male"
Set-TextValue "E4" "2001-04-16"
Set-TextValue "F4" "AnachalIdukki"
$ws.Range("G4").Value = 123456
